# Update dummy transaction ids in column C (rows 2 and 3):
# "test-id-14" -> "EMP-ID-25"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "EMP-ID-25"
$ws.Range("C3").Value = "EMP-ID-25"
